# Insert a new data row above row 724 (pushes the existing rows 724-810
# down to 725-811, matching the new sheet dimension A1:R811), then
# populate the newly inserted row 724 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(724).Insert()

$ws.Range('A724').Value2 = 3
$ws.Range('B724').Value2 = 'Femacal de La Calera'
$ws.Range('C724').Value2 = 'Coquimbo'
$ws.Range('D724').Value2 = 44918
$ws.Range('E724').Value2 = 5
$ws.Range('F724').Value2 = 100112006
$ws.Range('G724').Value2 = 'Repollo'
$ws.Range('H724').Value2 = 'Crespo record'
$ws.Range('I724').Value2 = 'Primera'
$ws.Range('J724').Value2 = 3150
$ws.Range('K724').Value2 = 1200
$ws.Range('L724').Value2 = 1300
$ws.Range('M724').Value2 = 1252
$ws.Range('N724').Value2 = '$/unidad'
$ws.Range('O724').Value2 = 'Provincia de Quillota'
$ws.Range('P724').Value2 = 1252
$ws.Range('Q724').Value2 = 1
$ws.Range('R724').Value2 = 'Hortaliza'
